$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '50.943.56'
$ws.Range("E2").Value = '  -1.11%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.926.75'
$ws.Range("E3").Value = '  -1.86%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '372.69'
$ws.Range("E5").Value = '  -1.92%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.01'
$ws.Range("E6").Value = '  -4.27%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.534'
$ws.Range("E7").Value = '  -1.38%  '
$ws.Range("E8").Value = '  +0.18%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.580'
$ws.Range("E9").Value = '  -2.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.92'
$ws.Range("E10").Value = '  -4.01%  '
$ws.Range("E11").Value = '  -0.82%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0841'
$ws.Range("E12").Value = '  -0.76%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.394.46'
$ws.Range("E13").Value = '  -1.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.92'
$ws.Range("E14").Value = '  -3.03%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.41'
$ws.Range("E15").Value = '  -1.90%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '11.70'
$ws.Range("E16").Value = '  +57.46%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.924.85'
$ws.Range("E17").Value = '  -1.72%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.963'
$ws.Range("E18").Value = '  -0.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '50.935.40'
$ws.Range("E19").Value = '  -1.06%  '
$ws.Range("E20").Value = '  -6.68%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.41'
$ws.Range("E21").Value = '  -4.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0952'
$ws.Range("E22").Value = '  -1.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '264.40'
$ws.Range("E23").Value = '  +0.85%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.34'
$ws.Range("E24").Value = '  -1.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.10'
$ws.Range("E25").Value = '  +9.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.99'
$ws.Range("E26").Value = '  +3.37%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.38'
$ws.Range("E27").Value = '  -3.08%  '
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.165'
$ws.Range("E29").Value = '  -4.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '25.45'
$ws.Range("E30").Value = '  -1.92%  '
$ws.Range("E31").Value = '  -2.31%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.97'
$ws.Range("E32").Value = '  +0.41%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '50.45'
$ws.Range("E33").Value = '  -1.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.03'
$ws.Range("E34").Value = '  -2.90%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0441'
$ws.Range("E35").Value = '  -1.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '32.79'
$ws.Range("E36").Value = '  -7.47%  '
$ws.Range("E37").Value = '  -0.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.20'
$ws.Range("E38").Value = '  +4.14%  '
$ws.Range("E39").Value = '  -1.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '16.26'
$ws.Range("E40").Value = '  -5.96%  '
$ws.Range("E41").Value = '  -5.95%  '
$ws.Range("E42").Value = '  -4.40%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '119.62'
$ws.Range("E43").Value = '  -4.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.13'
$ws.Range("E44").Value = '  -3.66%  '
$ws.Range("E45").Value = '  -0.41%  '
$ws.Range("E46").Value = '  -6.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.28'
$ws.Range("E47").Value = '  +1.36%  '
$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.30'
$ws.Range("E48").Value = '  -2.60%  '
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.995.40'
$ws.Range("E49").Value = '  -2.54%  '
$ws.Range("E50").Value = '  -4.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.30'
$ws.Range("E51").Value = '  +0.99%  '
